$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.939.59'
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").Value = '2.360.97'
$ws.Range("E3").Value = '  +1.43%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''302.99'
$ws.Range("E5").Value = '  +0.38%  '
$ws.Range("D6").Value = '''95.52'
$ws.Range("E6").Value = '  -0.39%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E8").Value = '  -0.54%  '
$ws.Range("D9").Value = '''0.482'
$ws.Range("E9").Value = '  -2.68%  '
$ws.Range("D10").Value = '''34.07'
$ws.Range("E10").Value = '  -0.78%  '
$ws.Range("D11").Value = '''0.124'
$ws.Range("E11").Value = '  +3.00%  '
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").Value = '''18.44'
$ws.Range("E13").Value = '  -3.27%  '
$ws.Range("E14").Value = '  -0.67%  '
$ws.Range("D15").Value = '2.727.21'
$ws.Range("E15").Value = '  +1.30%  '
$ws.Range("D16").Value = '2.348.66'
$ws.Range("E16").Value = '  +0.40%  '
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("D18").Value = '42.917.15'
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").Value = '''11.90'
$ws.Range("E19").Value = '  -2.99%  '
$ws.Range("D20").Value = '''6.25'
$ws.Range("E20").Value = '  +1.24%  '
$ws.Range("D21").Value = '0.0₃0884'
$ws.Range("E21").Value = '  -0.96%  '
$ws.Range("D22").Value = '''67.96'
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Value = '''234.94'
$ws.Range("E23").Value = '  -0.67%  '
$ws.Range("E24").Value = '  -4.25%  '
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").Value = '''2.43'
$ws.Range("E25").Value = '  +0.99%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").Value = '''24.46'
$ws.Range("E27").Value = '  -0.90%  '
$ws.Range("E28").Value = '  +15.28%  '
$ws.Range("D29").Value = '''9.30'
$ws.Range("E29").Value = '  +1.82%  '
$ws.Range("E30").Value = '  -0.84%  '
$ws.Range("D31").Value = '''1.00'
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("D32").Value = '''4.99'
$ws.Range("E32").Value = '  -0.36%  '
$ws.Range("D33").Value = '''17.50'
$ws.Range("E33").Value = '  -2.32%  '
$ws.Range("D34").Value = '''0.0715'
$ws.Range("E34").Value = '  +1.79%  '
$ws.Range("D35").Value = '''128.46'
$ws.Range("E35").Value = '  -11.21%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '''0.104'
$ws.Range("E36").Value = '  +3.13%  '
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").Value = '''1.83'
$ws.Range("E37").Value = '  +1.40%  '
$ws.Range("D38").Value = '''4.31'
$ws.Range("E38").Value = '  -2.34%  '
$ws.Range("D39").Value = '''2.82'
$ws.Range("E39").Value = '  +2.92%  '
$ws.Range("E40").Value = '  -1.36%  '
$ws.Range("E41").Value = '  -1.08%  '
$ws.Range("D42").Value = '''21.23'
$ws.Range("E42").Value = '  -4.40%  '
$ws.Range("D43").Value = '1.928.18'
$ws.Range("E43").Value = '  -0.31%  '
$ws.Range("D44").Value = '''0.0277'
$ws.Range("E44").Value = '  -0.44%  '
$ws.Range("E45").Value = '  +3.29%  '
$ws.Range("D46").Value = '''9.22'
$ws.Range("E46").Value = '  -9.00%  '
$ws.Range("E47").Value = '  -2.06%  '
$ws.Range("D48").Value = '2.588.46'
$ws.Range("E49").Value = '  +1.09%  '
$ws.Range("D50").Value = '''71.45'
$ws.Range("E50").Value = '  -2.35%  '
$ws.Range("E51").Value = '  +0.95%  '
